$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Format the Price column as Text first so numeric-looking values like
# "0.998" are written as literal strings (matching the source data),
# rather than being auto-converted to numbers by Excel.
$ws.Range('D2:D51').NumberFormat = '@'

$ws.Range('D2').Value = '64.965.75'
$ws.Range('E2').Value = '  +5.39%  '
$ws.Range('D3').Value = '2.980.53'
$ws.Range('E3').Value = '  +3.15%  '
$ws.Range('D4').Value = '0.998'
$ws.Range('E4').Value = '  -0.38%  '
$ws.Range('D5').Value = '580.02'
$ws.Range('E5').Value = '  +2.01%  '
$ws.Range('D6').Value = '152.68'
$ws.Range('E6').Value = '  +7.51%  '
$ws.Range('D7').Value = '0.999'
$ws.Range('E7').Value = '  -0.23%  '
$ws.Range('D8').Value = '2.978.97'
$ws.Range('E8').Value = '  +3.14%  '
$ws.Range('D9').Value = '0.513'
$ws.Range('E9').Value = '  +1.42%  '
$ws.Range('E10').Value = '  +4.92%  '
$ws.Range('D11').Value = '0.151'
$ws.Range('E11').Value = '  +3.39%  '
$ws.Range('D12').Value = '0.447'
$ws.Range('E12').Value = '  +3.10%  '
$ws.Range('D13').Value = '0.0000238'
$ws.Range('E13').Value = '  +3.11%  '
$ws.Range('D14').Value = '34.22'
$ws.Range('E14').Value = '  +7.47%  '
$ws.Range('E15').Value = '  +0.80%  '
$ws.Range('D16').Value = '64.824.14'
$ws.Range('E16').Value = '  +5.12%  '
$ws.Range('D17').Value = '3.473.21'
$ws.Range('E17').Value = '  +3.07%  '
$ws.Range('E18').Value = '  +3.99%  '
$ws.Range('D19').Value = '2.986.82'
$ws.Range('E19').Value = '  +3.57%  '
$ws.Range('D20').Value = '447.73'
$ws.Range('E20').Value = '  +3.51%  '
$ws.Range('D21').Value = '13.68'
$ws.Range('E21').Value = '  +3.89%  '
$ws.Range('D22').Value = '0.678'
$ws.Range('E22').Value = '  +3.63%  '
$ws.Range('D23').Value = '7.27'
$ws.Range('E23').Value = '  +5.72%  '
$ws.Range('D24').Value = '80.86'
$ws.Range('E24').Value = '  +1.67%  '
$ws.Range('D25').Value = '10.69'
$ws.Range('E25').Value = '  +5.04%  '
$ws.Range('E26').Value = '  +3.73%  '
$ws.Range('D27').Value = '2.19'
$ws.Range('E27').Value = '  +8.31%  '
$ws.Range('E28').Value = '  +0.03%  '
$ws.Range('E29').Value = '  +15.38%  '
$ws.Range('D30').Value = '7.75'
$ws.Range('E30').Value = '  +10.75%  '
$ws.Range('D31').Value = '0.0000102'
$ws.Range('E31').Value = '  +0.34%  '
$ws.Range('E32').Value = '  +2.94%  '
$ws.Range('E33').Value = '  +3.52%  '
$ws.Range('D34').Value = '26.61'
$ws.Range('E34').Value = '  +4.29%  '
$ws.Range('D35').Value = '0.998'
$ws.Range('E35').Value = '  -0.36%  '
$ws.Range('D36').Value = '0.982'
$ws.Range('E36').Value = '  +2.36%  '
$ws.Range('D37').Value = '5.65'
$ws.Range('E37').Value = '  +4.50%  '
$ws.Range('D38').Value = '2.09'
$ws.Range('E38').Value = '  +8.18%  '
$ws.Range('D39').Value = '48.92'
$ws.Range('E39').Value = '  -0.09%  '
$ws.Range('D40').Value = '2.88'
$ws.Range('E40').Value = '  +3.25%  '
$ws.Range('D41').Value = '43.71'
$ws.Range('E41').Value = '  +11.52%  '
$ws.Range('E42').Value = '  +3.61%  '
$ws.Range('D43').Value = '0.296'
$ws.Range('E43').Value = '  +11.06%  '
$ws.Range('D44').Value = '8.40'
$ws.Range('E44').Value = '  +2.05%  '
$ws.Range('D45').Value = '380.38'
$ws.Range('E45').Value = '  +12.69%  '
$ws.Range('D46').Value = '2.772.00'
$ws.Range('E46').Value = '  +2.83%  '
$ws.Range('E47').Value = '  +4.49%  '
$ws.Range('D48').Value = '133.98'
$ws.Range('E48').Value = '  +0.93%  '
$ws.Range('E49').Value = '  +0.00%  '
$ws.Range('E50').Value = '  +2.11%  '
$ws.Range('B51').Value = 'InjectiveProtocol'
$ws.Range('C51').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D51').Value = '22.84'
$ws.Range('E51').Value = '  +5.92%  '
